# Auto-generated edit script to update live market-price / profit columns (H:N)
# on multiple worksheets, per the scheduled-runner market data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1057.0834
$ws.Range("I2").Value = 298.5
$ws.Range("J2").Value = 1436.375
$ws.Range("K2").Value = 298.5
$ws.Range("L2").Value = 1436.375
$ws.Range("M2").Value = -185.5
$ws.Range("N2").Value = -1662.375
$ws.Range("H118").Value = 301.33334
$ws.Range("I118").Value = 197.5
$ws.Range("K118").Value = 592.5
$ws.Range("M118").Value = 1064.5
$ws.Range("H135").Value = 851
$ws.Range("I135").Value = 829.7143
$ws.Range("K135").Value = 7467.428699999999
$ws.Range("M135").Value = -4932.428699999999
$ws.Range("H137").Value = 4387.2104
$ws.Range("I137").Value = 4025.9
$ws.Range("K137").Value = 12077.7
$ws.Range("M137").Value = -9527.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -327
$ws.Range("H7").Value = 40000
$ws.Range("J7").Value = 40000
$ws.Range("L7").Value = 40000
$ws.Range("N7").Value = -40228
$ws.Range("H32").Value = 3995.1082
$ws.Range("I32").Value = 1406.5
$ws.Range("K32").Value = 1406.5
$ws.Range("M32").Value = -1119.5
$ws.Range("H61").Value = 3331
$ws.Range("I61").Value = 2996.75
$ws.Range("K61").Value = 2996.75
$ws.Range("M61").Value = -2784.75
$ws.Range("H62").Value = 45833
$ws.Range("J62").Value = 45833
$ws.Range("L62").Value = 45833
$ws.Range("N62").Value = -47081
$ws.Range("H65").Value = 45833
$ws.Range("J65").Value = 45833
$ws.Range("L65").Value = 137499
$ws.Range("N65").Value = -143739
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H135").Value = 29997.5
$ws.Range("J135").Value = 29997.5
$ws.Range("L135").Value = 29997.5
$ws.Range("N135").Value = -40137.5
$ws.Range("H136").Value = 3331
$ws.Range("I136").Value = 2996.75
$ws.Range("K136").Value = 8990.25
$ws.Range("M136").Value = -6440.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 30345.637
$ws.Range("J19").Value = 50000
$ws.Range("L19").Value = 50000
$ws.Range("N19").Value = -50346
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H82").Value = 17481.154
$ws.Range("H85").Value = 17481.154
$ws.Range("H99").Value = 1361.5
$ws.Range("I99").Value = 1415.3334
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 1415.3334
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = 82.66660000000002
$ws.Range("N99").Value = -4196

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 998.75
$ws.Range("J4").Value = 998.75
$ws.Range("L4").Value = 998.75
$ws.Range("N4").Value = -1222.75
$ws.Range("H132").Value = 1998
$ws.Range("I132").Value = 493
$ws.Range("K132").Value = 1479
$ws.Range("M132").Value = 1051
$ws.Range("H134").Value = 6189
$ws.Range("I134").Value = 5644.5713
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 16933.7139
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -14398.7139
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 460.4
$ws.Range("J2").Value = 1001
$ws.Range("L2").Value = 6006
$ws.Range("N2").Value = -6232
$ws.Range("H33").Value = 169.8
$ws.Range("J33").Value = 116.666664
$ws.Range("L33").Value = 699.999984
$ws.Range("N33").Value = -1265.999984
$ws.Range("H61").Value = 181.9
$ws.Range("J61").Value = 317.5
$ws.Range("L61").Value = 952.5
$ws.Range("N61").Value = -1382.5
$ws.Range("H122").Value = 993.3333
$ws.Range("J122").Value = 992.5
$ws.Range("L122").Value = 8932.5
$ws.Range("N122").Value = -13832.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 39226
$ws.Range("I63").Value = 38790
$ws.Range("J63").Value = 39444
$ws.Range("K63").Value = 38790
$ws.Range("L63").Value = 39444
$ws.Range("M63").Value = -38104
$ws.Range("N63").Value = -40816
$ws.Range("H66").Value = 39226
$ws.Range("I66").Value = 38790
$ws.Range("J66").Value = 39444
$ws.Range("K66").Value = 116370
$ws.Range("L66").Value = 118332
$ws.Range("M66").Value = -112938
$ws.Range("N66").Value = -125196
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22298
$ws.Range("I7").Value = 19378.727
$ws.Range("J7").Value = 27650
$ws.Range("K7").Value = 19378.727
$ws.Range("L7").Value = 27650
$ws.Range("M7").Value = -19266.727
$ws.Range("N7").Value = -27874
$ws.Range("H22").Value = 716.76
$ws.Range("I22").Value = 657.34784
$ws.Range("K22").Value = 657.34784
$ws.Range("M22").Value = -362.34784
$ws.Range("H27").Value = 716.76
$ws.Range("I27").Value = 657.34784
$ws.Range("K27").Value = 657.34784
$ws.Range("M27").Value = -550.34784
$ws.Range("H40").Value = 4667.3335
$ws.Range("I40").Value = 3874.75
$ws.Range("K40").Value = 3874.75
$ws.Range("M40").Value = -3738.75
$ws.Range("H61").Value = 3948.875
$ws.Range("I61").Value = 3598.6667
$ws.Range("J61").Value = 4999.5
$ws.Range("K61").Value = 3598.6667
$ws.Range("L61").Value = 4999.5
$ws.Range("M61").Value = -3396.6667
$ws.Range("N61").Value = -5403.5
$ws.Range("H62").Value = 55000
$ws.Range("J62").Value = 55000
$ws.Range("L62").Value = 55000
$ws.Range("N62").Value = -56248
$ws.Range("H65").Value = 55000
$ws.Range("J65").Value = 55000
$ws.Range("L65").Value = 165000
$ws.Range("N65").Value = -171240
$ws.Range("H93").Value = 3266.6667
$ws.Range("I93").Value = 4400
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 4400
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -3152
$ws.Range("N93").Value = -3496
$ws.Range("H113").Value = 3948.875
$ws.Range("I113").Value = 3598.6667
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 3598.6667
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = -1428.6667
$ws.Range("N113").Value = -9339.5
$ws.Range("H122").Value = 1045
$ws.Range("I122").Value = 1045
$ws.Range("K122").Value = 3135
$ws.Range("M122").Value = -685
$ws.Range("H126").Value = 22298
$ws.Range("I126").Value = 19378.727
$ws.Range("J126").Value = 27650
$ws.Range("K126").Value = 58136.181
$ws.Range("L126").Value = 82950
$ws.Range("M126").Value = -55666.181
$ws.Range("N126").Value = -87890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H86").Value = 130000
$ws.Range("J86").Value = 130000
$ws.Range("L86").Value = 130000
$ws.Range("N86").Value = -132246
$ws.Range("H89").Value = 130000
$ws.Range("J89").Value = 130000
$ws.Range("L89").Value = 650000
$ws.Range("N89").Value = -661232
